$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column P (rows 3-10) into the new column Q so the
# new cells inherit the same styles (borders/fonts) as the rest of the
# table before we fill in the 2023 values.
$ws.Range("P3:P10").Copy()
$ws.Range("Q3:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "2023" column header
$ws.Range("Q4").Value = 2023

# New data values for the 2023 column
$ws.Range("Q6").Value = 1209
$ws.Range("Q7").Value = "-"
$ws.Range("Q8").Value = 373
$ws.Range("Q9").Value = 115
$ws.Range("Q10").Value = 781

# Row 5 grew slightly taller to fit the extra column
$ws.Rows.Item(5).RowHeight = 27

$ws.Range("A1").Select()
